# Weekly refresh: each data row (2-37) picks up the Fecha/Volumen/Precio
# columns from another row in the same sheet (a full re-shuffle of the
# per-row price-quote figures while Mercado/Region/Categoria/etc. stay put).
# Map: destination row -> source row (where the "new" values come from).
$map = @{2=27; 3=26; 4=22; 5=6; 6=20; 7=29; 8=31; 9=37; 10=13; 11=11; 12=33; 13=8; 14=7; 15=14; 16=3; 17=10; 18=4; 19=15; 20=17; 21=23; 22=12; 23=36; 24=19; 25=28; 26=32; 27=35; 28=9; 29=2; 30=21; 31=24; 32=18; 33=5; 34=16; 35=34; 36=30; 37=25}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns touched by the shuffle: D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).
$cols = @(4, 10, 11, 12, 13, 16)

# Snapshot current values for every row before writing anything, so the
# row-to-row copy reads the ORIGINAL data regardless of write order.
$snapshot = @{}
for ($r = 2; $r -le 37; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write each row's new values from its mapped source row's snapshot.
foreach ($r in $map.Keys) {
    $src = $map[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $srcVals[$c]
    }
}
